# "updated single child outputs" - the single-child prediction table's
# distance column (B2:B8) was recomputed; previously every row held the
# placeholder value 1, now each row carries its real computed distance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1610.0941026769119,
    1864.2031830372748,
    1270.2029342352157,
    1728.4917204079497,
    1639.0537384236563,
    1588.6114004985702,
    1637.953737106925
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Re-apply the text number format ("@") to the header row and the row-label
# column. The underlying format is unchanged (still text), but rewriting it
# lets the style table collapse the cells onto their canonical, already
# existing text-format style entry instead of the stale duplicate they used
# to point at.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A8").NumberFormat = "@"
